$wb = $excel.ActiveWorkbook

# --- ALC (sheet index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H62").Value = 932
$ws.Range("I62").Value = 500
$ws.Range("J62").Value = 1004
$ws.Range("K62").Value = 500
$ws.Range("L62").Value = 1004
$ws.Range("M62").Value = 124
$ws.Range("N62").Value = -2252
$ws.Range("H65").Value = 932
$ws.Range("I65").Value = 500
$ws.Range("J65").Value = 1004
$ws.Range("K65").Value = 2500
$ws.Range("L65").Value = 5020
$ws.Range("M65").Value = 620
$ws.Range("N65").Value = -11260
$ws.Range("H137").Value = 1568.25
$ws.Range("I137").Value = 1174.3572
$ws.Range("J137").Value = 1962.1428
$ws.Range("K137").Value = 3523.0716
$ws.Range("L137").Value = 5886.428400000001
$ws.Range("M137").Value = -973.0715999999998
$ws.Range("N137").Value = -10986.4284

# --- ARM (sheet index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1569.3889
$ws.Range("I2").Value = 1604.0588
$ws.Range("K2").Value = 1604.0588
$ws.Range("M2").Value = -1491.0588
$ws.Range("H32").Value = 4118.5537
$ws.Range("I32").Value = 3672.5107
$ws.Range("J32").Value = 5283.222
$ws.Range("K32").Value = 3672.5107
$ws.Range("L32").Value = 5283.222
$ws.Range("M32").Value = -3385.5107
$ws.Range("N32").Value = -5857.222
$ws.Range("H45").Value = 11095.2
$ws.Range("I45").Value = 13426.167
$ws.Range("J45").Value = 1771.3334
$ws.Range("K45").Value = 13426.167
$ws.Range("L45").Value = 1771.3334
$ws.Range("M45").Value = -13049.167
$ws.Range("N45").Value = -2525.3334
$ws.Range("H97").Value = 750.8
$ws.Range("I97").Value = 757.55554
$ws.Range("J97").Value = 690
$ws.Range("K97").Value = 757.55554
$ws.Range("L97").Value = 690
$ws.Range("M97").Value = -261.55554
$ws.Range("N97").Value = -1682
$ws.Range("H102").Value = 5306582
$ws.Range("I102").Value = 5306582
$ws.Range("K102").Value = 5306582
$ws.Range("M102").Value = -5304960
$ws.Range("H116").Value = 1569.3889
$ws.Range("I116").Value = 1604.0588
$ws.Range("K116").Value = 1604.0588
$ws.Range("M116").Value = 689.9412
$ws.Range("H122").Value = 2851231.8
$ws.Range("I122").Value = 6411896.5
$ws.Range("K122").Value = 19235689.5
$ws.Range("M122").Value = -19233239.5

# --- BSM (sheet index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1569.3889
$ws.Range("I3").Value = 1604.0588
$ws.Range("K3").Value = 1604.0588
$ws.Range("M3").Value = -1490.0588
$ws.Range("H105").Value = 13254.333
$ws.Range("I105").Value = 21670.9
$ws.Range("J105").Value = 2733.625
$ws.Range("K105").Value = 21670.9
$ws.Range("L105").Value = 2733.625
$ws.Range("M105").Value = -19923.9
$ws.Range("N105").Value = -6227.625
$ws.Range("H107").Value = 1476.4
$ws.Range("I107").Value = 1470.5
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1470.5
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 449.5
$ws.Range("N107").Value = -5340

# --- CRP (sheet index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 466.66666
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -100
$ws.Range("H31").Value = 3055.3833
$ws.Range("I31").Value = 2623.3928
$ws.Range("K31").Value = 2623.3928
$ws.Range("M31").Value = -2328.3928
$ws.Range("H34").Value = 3055.3833
$ws.Range("I34").Value = 2623.3928
$ws.Range("K34").Value = 2623.3928
$ws.Range("M34").Value = -2421.3928
$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 20000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 20000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -20812
$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 20000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 20000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -22808
$ws.Range("H99").Value = 20259.934
$ws.Range("I99").Value = 17285.715
$ws.Range("J99").Value = 22862.375
$ws.Range("K99").Value = 17285.715
$ws.Range("L99").Value = 22862.375
$ws.Range("M99").Value = -15787.715
$ws.Range("N99").Value = -25858.375
$ws.Range("H126").Value = 20259.934
$ws.Range("I126").Value = 17285.715
$ws.Range("J126").Value = 22862.375
$ws.Range("K126").Value = 51857.145
$ws.Range("L126").Value = 68587.125
$ws.Range("M126").Value = -49387.145
$ws.Range("N126").Value = -73527.125
$ws.Range("H132").Value = 2229.3667
$ws.Range("I132").Value = 1884.5555
$ws.Range("J132").Value = 5332.6665
$ws.Range("K132").Value = 5653.666499999999
$ws.Range("L132").Value = 15997.9995
$ws.Range("M132").Value = -3123.666499999999
$ws.Range("N132").Value = -21057.9995
$ws.Range("H134").Value = 3527.0435
$ws.Range("I134").Value = 3527.0435
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10581.1305
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8046.130500000001
$ws.Range("N134").ClearContents()

# --- CUL (sheet index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H74").Value = 17200
$ws.Range("J74").Value = 17200
$ws.Range("L74").Value = 51600
$ws.Range("N74").Value = -53722
$ws.Range("H77").Value = 17200
$ws.Range("J77").Value = 17200
$ws.Range("L77").Value = 154800
$ws.Range("N77").Value = -165408
$ws.Range("H131").Value = 18645096
$ws.Range("I131").Value = 5882964
$ws.Range("J131").Value = 23810720
$ws.Range("K131").Value = 17648892
$ws.Range("L131").Value = 71432160
$ws.Range("M131").Value = -17643852
$ws.Range("N131").Value = -71442240

# --- GSM (sheet index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 132.31818
$ws.Range("I2").Value = 172
$ws.Range("K2").Value = 172
$ws.Range("M2").Value = -59
$ws.Range("H126").Value = 11276.762
$ws.Range("I126").Value = 13459.529
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 40378.587
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -37908.587
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 3529.0754
$ws.Range("I132").Value = 3397.625
$ws.Range("J132").Value = 3637.862
$ws.Range("K132").Value = 10192.875
$ws.Range("L132").Value = 10913.586
$ws.Range("M132").Value = -7662.875
$ws.Range("N132").Value = -15973.586

# --- LTW (sheet index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H59").Value = 27398
$ws.Range("J59").Value = 27398
$ws.Range("L59").Value = 27398
$ws.Range("N59").Value = -28706
$ws.Range("H122").Value = 6788170
$ws.Range("I122").Value = 7145304.5
$ws.Range("J122").Value = 5002500
$ws.Range("K122").Value = 21435913.5
$ws.Range("L122").Value = 15007500
$ws.Range("M122").Value = -21433463.5
$ws.Range("N122").Value = -15012400
$ws.Range("H138").Value = 59789.5
$ws.Range("J138").Value = 59789.5
$ws.Range("L138").Value = 59789.5
$ws.Range("N138").Value = -70069.5

# --- WVR (sheet index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H80").Value = 39301
$ws.Range("J80").Value = 39301
$ws.Range("L80").Value = 39301
$ws.Range("N80").Value = -41297
$ws.Range("H83").Value = 39301
$ws.Range("J83").Value = 39301
$ws.Range("L83").Value = 117903
$ws.Range("N83").Value = -127887
